$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "DPF start"
$ws.Range("D13").Interior.Color = 65535
$ws.Range("F14").Value = "DPF lamp"
$ws.Range("F14").Interior.Color = 65535
$ws.Range("D6").Value = "H window button"
$ws.Range("B17").Value = "heated window"
$ws.Range("B19").Value = "heated window"
$ws.Range("F17").Value = "idle solenoid"
$ws.Range("D16").Value = "turbo solenoid"

$ws.Range("D16").Select()
